$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write the "Temps joue" (G) values first, in row order,
# so new shared strings are appended in the same order as the target file.
$ws.Cells.Item(567,7).Value = "01:37:55"
$ws.Cells.Item(568,7).Value = "01:36:57"
$ws.Cells.Item(569,7).Value = "01:38:17"
$ws.Cells.Item(570,7).Value = "01:38:25"
$ws.Cells.Item(571,7).Value = "01:36:56"
$ws.Cells.Item(572,7).Value = "00:47:56"
$ws.Cells.Item(573,7).Value = "01:36:57"
$ws.Cells.Item(574,7).Value = "00:59:51"
$ws.Cells.Item(575,7).Value = "01:19:04"
$ws.Cells.Item(576,7).Value = "00:18:22"
$ws.Cells.Item(577,7).Value = "00:59:36"

# Step 2: write the match/session name (A) - this becomes the last new shared string.
$ws.Cells.Item(567,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(568,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(569,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(570,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(571,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(572,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(573,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(574,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(575,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(576,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"
$ws.Cells.Item(577,1).Value = "CDF T4 VS Misérieux Trévoux (R1)"

# Step 3: date (B), with the same date-number-format style as the existing rows.
$ws.Range("B566").Copy()
$ws.Range("B567:B577").PasteSpecial(-4122)
$ws.Cells.Item(567,2).Value = 45927
$ws.Cells.Item(568,2).Value = 45927
$ws.Cells.Item(569,2).Value = 45927
$ws.Cells.Item(570,2).Value = 45927
$ws.Cells.Item(571,2).Value = 45927
$ws.Cells.Item(572,2).Value = 45927
$ws.Cells.Item(573,2).Value = 45927
$ws.Cells.Item(574,2).Value = 45927
$ws.Cells.Item(575,2).Value = 45927
$ws.Cells.Item(576,2).Value = 45927
$ws.Cells.Item(577,2).Value = 45927

# Step 4: remaining text columns (C, E, F).
$ws.Cells.Item(567,3).Value = "Global"
$ws.Cells.Item(567,5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(567,6).Value = "center midfield"
$ws.Cells.Item(568,3).Value = "Global"
$ws.Cells.Item(568,5).Value = "Kamal Bafounta"
$ws.Cells.Item(568,6).Value = "center midfield"
$ws.Cells.Item(569,3).Value = "Global"
$ws.Cells.Item(569,5).Value = "Naim Ighbane"
$ws.Cells.Item(569,6).Value = "center back"
$ws.Cells.Item(570,3).Value = "Global"
$ws.Cells.Item(570,5).Value = "Naim Dhib"
$ws.Cells.Item(570,6).Value = "center midfield"
$ws.Cells.Item(571,3).Value = "Global"
$ws.Cells.Item(571,5).Value = "Amir Etien"
$ws.Cells.Item(571,6).Value = "right forward"
$ws.Cells.Item(572,3).Value = "Global"
$ws.Cells.Item(572,5).Value = "Omar Benyounes"
$ws.Cells.Item(572,6).Value = "center midfield"
$ws.Cells.Item(573,3).Value = "Global"
$ws.Cells.Item(573,5).Value = "Yoan Zouma"
$ws.Cells.Item(573,6).Value = "center back"
$ws.Cells.Item(574,3).Value = "Global"
$ws.Cells.Item(574,5).Value = "Karim Belmahi"
$ws.Cells.Item(574,6).Value = "left forward"
$ws.Cells.Item(575,3).Value = "Global"
$ws.Cells.Item(575,5).Value = "Malik Boussaid"
$ws.Cells.Item(575,6).Value = "right back"
$ws.Cells.Item(576,3).Value = "Global"
$ws.Cells.Item(576,5).Value = "Hedi Nasri"
$ws.Cells.Item(576,6).Value = "right back"
$ws.Cells.Item(577,3).Value = "Global"
$ws.Cells.Item(577,5).Value = "Emmanuel Valey"
$ws.Cells.Item(577,6).Value = "left forward"

# Step 5: numeric stat columns H..V.
$ws.Cells.Item(567,8).Value = 11.8
$ws.Cells.Item(567,9).Value = 2.43
$ws.Cells.Item(567,10).Value = 9.35
$ws.Cells.Item(567,11).Value = 1.73
$ws.Cells.Item(567,12).Value = 0.59
$ws.Cells.Item(567,13).Value = 0.13
$ws.Cells.Item(567,14).Value = 0.0
$ws.Cells.Item(567,15).Value = 9.0
$ws.Cells.Item(567,16).Value = 7.15
$ws.Cells.Item(567,17).Value = 29.74
$ws.Cells.Item(567,18).Value = 4.58
$ws.Cells.Item(567,19).Value = 47.0
$ws.Cells.Item(567,20).Value = 7.0
$ws.Cells.Item(567,21).Value = 23.0
$ws.Cells.Item(567,22).Value = 11.0
$ws.Cells.Item(568,8).Value = 11.12
$ws.Cells.Item(568,9).Value = 2.22
$ws.Cells.Item(568,10).Value = 8.88
$ws.Cells.Item(568,11).Value = 1.71
$ws.Cells.Item(568,12).Value = 0.42
$ws.Cells.Item(568,13).Value = 0.12
$ws.Cells.Item(568,14).Value = 0.0
$ws.Cells.Item(568,15).Value = 5.0
$ws.Cells.Item(568,16).Value = 6.91
$ws.Cells.Item(568,17).Value = 29.54
$ws.Cells.Item(568,18).Value = 4.34
$ws.Cells.Item(568,19).Value = 30.0
$ws.Cells.Item(568,20).Value = 2.0
$ws.Cells.Item(568,21).Value = 33.0
$ws.Cells.Item(568,22).Value = 4.0
$ws.Cells.Item(569,8).Value = 9.3
$ws.Cells.Item(569,9).Value = 1.05
$ws.Cells.Item(569,10).Value = 8.24
$ws.Cells.Item(569,11).Value = 0.66
$ws.Cells.Item(569,12).Value = 0.3
$ws.Cells.Item(569,13).Value = 0.1
$ws.Cells.Item(569,14).Value = 0.0
$ws.Cells.Item(569,15).Value = 5.0
$ws.Cells.Item(569,16).Value = 5.57
$ws.Cells.Item(569,17).Value = 29.79
$ws.Cells.Item(569,18).Value = 4.31
$ws.Cells.Item(569,19).Value = 23.0
$ws.Cells.Item(569,20).Value = 2.0
$ws.Cells.Item(569,21).Value = 19.0
$ws.Cells.Item(569,22).Value = 8.0
$ws.Cells.Item(570,8).Value = 9.71
$ws.Cells.Item(570,9).Value = 1.77
$ws.Cells.Item(570,10).Value = 7.92
$ws.Cells.Item(570,11).Value = 1.05
$ws.Cells.Item(570,12).Value = 0.55
$ws.Cells.Item(570,13).Value = 0.2
$ws.Cells.Item(570,14).Value = 0.0
$ws.Cells.Item(570,15).Value = 16.0
$ws.Cells.Item(570,16).Value = 5.92
$ws.Cells.Item(570,17).Value = 29.94
$ws.Cells.Item(570,18).Value = 4.95
$ws.Cells.Item(570,19).Value = 32.0
$ws.Cells.Item(570,20).Value = 3.0
$ws.Cells.Item(570,21).Value = 33.0
$ws.Cells.Item(570,22).Value = 20.0
$ws.Cells.Item(571,8).Value = 8.78
$ws.Cells.Item(571,9).Value = 1.95
$ws.Cells.Item(571,10).Value = 6.81
$ws.Cells.Item(571,11).Value = 0.96
$ws.Cells.Item(571,12).Value = 0.56
$ws.Cells.Item(571,13).Value = 0.32
$ws.Cells.Item(571,14).Value = 0.13
$ws.Cells.Item(571,15).Value = 20.0
$ws.Cells.Item(571,16).Value = 5.35
$ws.Cells.Item(571,17).Value = 34.73
$ws.Cells.Item(571,18).Value = 4.91
$ws.Cells.Item(571,19).Value = 46.0
$ws.Cells.Item(571,20).Value = 12.0
$ws.Cells.Item(571,21).Value = 22.0
$ws.Cells.Item(571,22).Value = 17.0
$ws.Cells.Item(572,8).Value = 6.02
$ws.Cells.Item(572,9).Value = 1.43
$ws.Cells.Item(572,10).Value = 4.57
$ws.Cells.Item(572,11).Value = 0.92
$ws.Cells.Item(572,12).Value = 0.36
$ws.Cells.Item(572,13).Value = 0.13
$ws.Cells.Item(572,14).Value = 0.03
$ws.Cells.Item(572,15).Value = 8.0
$ws.Cells.Item(572,16).Value = 7.5
$ws.Cells.Item(572,17).Value = 32.07
$ws.Cells.Item(572,18).Value = 4.69
$ws.Cells.Item(572,19).Value = 24.0
$ws.Cells.Item(572,20).Value = 5.0
$ws.Cells.Item(572,21).Value = 25.0
$ws.Cells.Item(572,22).Value = 13.0
$ws.Cells.Item(573,8).Value = 8.94
$ws.Cells.Item(573,9).Value = 1.13
$ws.Cells.Item(573,10).Value = 7.79
$ws.Cells.Item(573,11).Value = 0.8
$ws.Cells.Item(573,12).Value = 0.22
$ws.Cells.Item(573,13).Value = 0.12
$ws.Cells.Item(573,14).Value = 0.0
$ws.Cells.Item(573,15).Value = 6.0
$ws.Cells.Item(573,16).Value = 5.49
$ws.Cells.Item(573,17).Value = 29.51
$ws.Cells.Item(573,18).Value = 4.66
$ws.Cells.Item(573,19).Value = 30.0
$ws.Cells.Item(573,20).Value = 3.0
$ws.Cells.Item(573,21).Value = 24.0
$ws.Cells.Item(573,22).Value = 10.0
$ws.Cells.Item(574,8).Value = 6.68
$ws.Cells.Item(574,9).Value = 1.37
$ws.Cells.Item(574,10).Value = 5.29
$ws.Cells.Item(574,11).Value = 0.85
$ws.Cells.Item(574,12).Value = 0.38
$ws.Cells.Item(574,13).Value = 0.13
$ws.Cells.Item(574,14).Value = 0.02
$ws.Cells.Item(574,15).Value = 8.0
$ws.Cells.Item(574,16).Value = 6.66
$ws.Cells.Item(574,17).Value = 32.06
$ws.Cells.Item(574,18).Value = 4.9
$ws.Cells.Item(574,19).Value = 36.0
$ws.Cells.Item(574,20).Value = 8.0
$ws.Cells.Item(574,21).Value = 21.0
$ws.Cells.Item(574,22).Value = 14.0
$ws.Cells.Item(575,8).Value = 9.03
$ws.Cells.Item(575,9).Value = 1.97
$ws.Cells.Item(575,10).Value = 7.04
$ws.Cells.Item(575,11).Value = 1.22
$ws.Cells.Item(575,12).Value = 0.63
$ws.Cells.Item(575,13).Value = 0.14
$ws.Cells.Item(575,14).Value = 0.0
$ws.Cells.Item(575,15).Value = 11.0
$ws.Cells.Item(575,16).Value = 6.71
$ws.Cells.Item(575,17).Value = 27.72
$ws.Cells.Item(575,18).Value = 4.42
$ws.Cells.Item(575,19).Value = 31.0
$ws.Cells.Item(575,20).Value = 8.0
$ws.Cells.Item(575,21).Value = 26.0
$ws.Cells.Item(575,22).Value = 9.0
$ws.Cells.Item(576,8).Value = 1.99
$ws.Cells.Item(576,9).Value = 0.38
$ws.Cells.Item(576,10).Value = 1.61
$ws.Cells.Item(576,11).Value = 0.19
$ws.Cells.Item(576,12).Value = 0.1
$ws.Cells.Item(576,13).Value = 0.07
$ws.Cells.Item(576,14).Value = 0.02
$ws.Cells.Item(576,15).Value = 3.0
$ws.Cells.Item(576,16).Value = 6.48
$ws.Cells.Item(576,17).Value = 31.39
$ws.Cells.Item(576,18).Value = 4.07
$ws.Cells.Item(576,19).Value = 8.0
$ws.Cells.Item(576,20).Value = 2.0
$ws.Cells.Item(576,21).Value = 4.0
$ws.Cells.Item(576,22).Value = 3.0
$ws.Cells.Item(577,8).Value = 7.2
$ws.Cells.Item(577,9).Value = 1.58
$ws.Cells.Item(577,10).Value = 5.59
$ws.Cells.Item(577,11).Value = 1.03
$ws.Cells.Item(577,12).Value = 0.41
$ws.Cells.Item(577,13).Value = 0.12
$ws.Cells.Item(577,14).Value = 0.04
$ws.Cells.Item(577,15).Value = 9.0
$ws.Cells.Item(577,16).Value = 7.22
$ws.Cells.Item(577,17).Value = 32.64
$ws.Cells.Item(577,18).Value = 4.88
$ws.Cells.Item(577,19).Value = 38.0
$ws.Cells.Item(577,20).Value = 6.0
$ws.Cells.Item(577,21).Value = 23.0
$ws.Cells.Item(577,22).Value = 14.0

# Step 6: restore selection to reflect the saved view state.
$ws.Range("D583").Select()
